$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.002.64'
$ws.Range("E2").Value = '  -2.96%  '
$ws.Range("D3").Value = '3.368.91'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.35'
$ws.Range("E5").Value = '  -1.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.97'
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.97'
$ws.Range("E9").Value = '  +0.57%  '
$ws.Range("E10").Value = '  -0.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.416'
$ws.Range("E11").Value = '  +2.28%  '
$ws.Range("D12").Value = '3.948.74'
$ws.Range("E12").Value = '  -2.29%  '
$ws.Range("E13").Value = '  +0.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.02'
$ws.Range("E14").Value = '  -1.06%  '
$ws.Range("D15").Value = '3.361.43'
$ws.Range("E15").Value = '  -2.31%  '
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("D17").Value = '61.088.13'
$ws.Range("E17").Value = '  -2.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.34'
$ws.Range("E18").Value = '  -1.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.48'
$ws.Range("E19").Value = '  -1.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.92'
$ws.Range("E20").Value = '  -2.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '375.30'
$ws.Range("E21").Value = '  -3.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.34'
$ws.Range("E22").Value = '  +0.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.561'
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = '3.504.76'
$ws.Range("E25").Value = '  -2.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000109'
$ws.Range("E26").Value = '  -5.43%  '
$ws.Range("E27").Value = '  -3.87%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.48'
$ws.Range("E28").Value = '  -2.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  -1.19%  '
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.73'
$ws.Range("E32").Value = '  -3.81%  '
$ws.Range("E33").Value = '  -1.76%  '
$ws.Range("E34").Value = '  -3.51%  '
$ws.Range("E35").Value = '  +1.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '170.37'
$ws.Range("E36").Value = '  +0.13%  '
$ws.Range("E37").Value = '  -3.71%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.83'
$ws.Range("E38").Value = '  -2.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '29.06'
$ws.Range("E39").Value = '  -8.91%  '
$ws.Range("D40").Value = '3.403.63'
$ws.Range("E40").Value = '  -2.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0756'
$ws.Range("E41").Value = '  -3.51%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.761'
$ws.Range("E42").Value = '  -3.84%  '
$ws.Range("E43").Value = '  -1.20%  '
$ws.Range("E44").Value = '  -2.66%  '
$ws.Range("E45").Value = '  -5.72%  '
$ws.Range("D46").Value = '2.492.53'
$ws.Range("E46").Value = '  -2.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.70'
$ws.Range("E47").Value = '  -2.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.61'
$ws.Range("E48").Value = '  -0.24%  '
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("E50").Value = '  -2.30%  '
$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.816'
$ws.Range("E51").Value = '  +0.11%  '
